$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (e.g. "312.06", "0.3400", "-0.01%") are preserved verbatim instead of
# being auto-converted to numbers/percentages by Excel's type inference.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "312.06"
$ws.Range("D3").Value = "37.73"
$ws.Range("E3").Value = "0.06%"
$ws.Range("D4").Value = "5.135"
$ws.Range("E4").Value = "0.61%"
$ws.Range("E5").Value = "0.54%"
$ws.Range("D6").Value = "4.411"
$ws.Range("E6").Value = "0.93%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.905"
$ws.Range("E7").Value = "-3.00%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "8.274"
$ws.Range("E8").Value = "-0.38%"
$ws.Range("D9").Value = "0.9261"
$ws.Range("E9").Value = "-0.15%"
$ws.Range("D10").Value = "0.1203"
$ws.Range("E10").Value = "-10.89%"
$ws.Range("E11").Value = "-1.73%"
$ws.Range("D12").Value = "0.09164"
$ws.Range("E12").Value = "2.34%"
$ws.Range("D13").Value = "0.03296"
$ws.Range("E13").Value = "-5.18%"
$ws.Range("D14").Value = "0.09627"
$ws.Range("E14").Value = "-0.84%"
$ws.Range("D15").Value = "0.001386"
$ws.Range("E15").Value = "-0.55%"
$ws.Range("D16").Value = "0.005859"
$ws.Range("E16").Value = "-2.51%"
$ws.Range("E17").Value = "-2.17%"
$ws.Range("D18").Value = "3.096"
$ws.Range("E18").Value = "-0.26%"
$ws.Range("D19").Value = "0.3400"
$ws.Range("E19").Value = "-1.89%"
$ws.Range("D20").Value = "5.279"
$ws.Range("E20").Value = "5.50%"
$ws.Range("D21").Value = "0.1272"
$ws.Range("E21").Value = "-1.72%"
$ws.Range("E22").Value = "3.01%"
$ws.Range("D24").Value = "0.04361"
$ws.Range("E24").Value = "0.31%"
$ws.Range("D25").Value = "0.001249"
$ws.Range("E25").Value = "1.94%"
$ws.Range("E26").Value = "-5.07%"
$ws.Range("D27").Value = "0.0001220"
$ws.Range("E27").Value = "-9.74%"
$ws.Range("D39").Value = "0.02122"
$ws.Range("E39").Value = "-7.09%"
$ws.Range("D40").Value = "0.05159"
$ws.Range("E40").Value = "1.87%"
$ws.Range("D41").Value = "0.007663"
$ws.Range("E41").Value = "0.70%"
$ws.Range("E42").Value = "-6.61%"
$ws.Range("E43").Value = "0.39%"
$ws.Range("E44").Value = "-1.10%"
$ws.Range("D45").Value = "0.008603"
$ws.Range("E45").Value = "-2.07%"
$ws.Range("D46").Value = "0.00006683"
$ws.Range("E46").Value = "-1.87%"
$ws.Range("E47").Value = "-0.08%"
$ws.Range("D49").Value = "0.002869"
$ws.Range("E49").Value = "-4.47%"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.08%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.08%"

# Restore default (un-styled) cell formatting so only the values changed.
$ws.Range("D2:E51").Style = "Normal"
